# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts) across the resume's bullet points.
#
# Each target run of text is split into three (or more) runs: the
# surrounding plain text stays unformatted, while the metric itself
# becomes Bold with font color #2C3E50 (decimal 5258796 == RGB(44,62,80)
# packed as 0x00BBGGRR, which is how Word's Font.Color OLE property
# expects values).

$d = $word.ActiveDocument
$highlightColor = 5258796   # 0x2C3E50 -> BGR-packed OLE color

function Highlight-Next($SearchRange, $Text, $ParagraphEnd) {
    $SearchRange.Find.Execute($Text, $true, $false, $false, $false, $false,
                               $true, 1, $false, "", 0)
    $SearchRange.Font.Bold = $true
    $SearchRange.Font.Color = $highlightColor
    return $d.Range($SearchRange.End, $ParagraphEnd)
}

# --- Paragraph 10: "Discovered systematic race coding errors ..." ---
$p = $d.Paragraphs.Item(10)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "23%" $pEnd
$rng = Highlight-Next $rng "64%" $pEnd

# --- Paragraph 12: "Utilized advanced sampling methods ..." ---
$p = $d.Paragraphs.Item(12)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "±4.2%" $pEnd
$rng = Highlight-Next $rng "±2.1%" $pEnd
$rng = Highlight-Next $rng "71%" $pEnd
$rng = Highlight-Next $rng "87%" $pEnd

# --- Paragraph 13: "Trigonometric algorithm for boundary estimation ..." ---
$p = $d.Paragraphs.Item(13)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "73.5%" $pEnd
$rng = Highlight-Next $rng "$4.7M" $pEnd

# --- Paragraph 14: "Built real-time FEC analysis systems ..." ---
$p = $d.Paragraphs.Item(14)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "$2" $pEnd

# --- Paragraph 50: "Algorithmic innovation: Pioneered trigonometric ..." ---
$p = $d.Paragraphs.Item(50)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "73.5%" $pEnd

# --- Paragraph 51: "$4.7M savings enabled nonprofit access" ---
$p = $d.Paragraphs.Item(51)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "$4.7M" $pEnd

# --- Paragraph 53: "178% accuracy improvement in racial classification algorithms" ---
$p = $d.Paragraphs.Item(53)
$pEnd = $p.Range.End
$rng = $d.Range($p.Range.Start, $pEnd)
$rng = Highlight-Next $rng "178%" $pEnd
